# Insert a new weekly record row above row 582 (this pushes the existing
# rows 582..615 down to 583..616, preserving all of their original values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(582).Insert()

$ws.Cells.Item(582, 1).Value = 9
$ws.Cells.Item(582, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(582, 3).Value = "Metropolitana"
$ws.Cells.Item(582, 4).Value = 44753
$ws.Cells.Item(582, 5).Value = 13
$ws.Cells.Item(582, 6).Value = 100112009
$ws.Cells.Item(582, 7).Value = "Acelga"
$ws.Cells.Item(582, 8).Value = "Sin especificar"
$ws.Cells.Item(582, 9).Value = "Primera"
$ws.Cells.Item(582, 10).Value = 52
$ws.Cells.Item(582, 11).Value = 20000
$ws.Cells.Item(582, 12).Value = 20000
$ws.Cells.Item(582, 13).Value = 20000
$ws.Cells.Item(582, 14).Value = "`$/docena de atados"
$ws.Cells.Item(582, 15).Value = "Región Metropolitana"
$ws.Cells.Item(582, 16).Value = 6667
$ws.Cells.Item(582, 17).Value = 3
$ws.Cells.Item(582, 18).Value = "Hortaliza"
